# Combine sentiment results: add a Row_Number column and drop the
# trailing SAP (2021-06-24) and three Adidas (2021-06-25) rows that
# are no longer part of the combined data frame.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the last 4 data rows (SAP 2021-06-24 + 3x Adidas 2021-06-25)
$ws.Rows("23:26").Delete()

# 2. Insert a new column before the NLTK_Sentiment column (currently G)
#    to hold the new Row_Number field.
$ws.Columns("G:G").Insert()

# 3. Header for the new column
$ws.Range("G1").Value = "Row_Number"

# 4. Sequential row numbers (1-based) for each remaining data row
for ($r = 2; $r -le 22; $r++) {
    $ws.Range("G$r").Value = $r - 1
}

# 5. Match the header formatting used by the rest of the header row
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
